$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Q0)
$ws.Range("B3").Value = 0.2450315797533715
$ws.Range("C3").Value = 1.195258092086972
$ws.Range("D3").Value = 3.811732247447177
$ws.Range("E3").Value = 1.952365807795039
$ws.Range("F3").Value = 1.943689115874897
$ws.Range("G3").Value = 144

# Row 4 (Q1)
$ws.Range("B4").Value = 0.3216218552325188
$ws.Range("C4").Value = 1.318186357591293
$ws.Range("D4").Value = 8.12545414019174
$ws.Range("E4").Value = 2.850518223094134
$ws.Range("F4").Value = 2.852766129210556
$ws.Range("G4").Value = 70
